$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "APA Hotel Asakusa Ekimae"
$ws.Range("C2").Value = "Scored 8.2"
$ws.Range("D2").Value = "Very Good`n5,474 reviews"
$ws.Range("E2").Value = "50 m from center"
$ws.Rows.Item(2).AutoFit()

$ws.Range("A3").Value = "APA Hotel Roppongi SIX"
$ws.Range("C3").Value = "Scored 7.9"
$ws.Range("D3").Value = "Good`n10,789 reviews"
$ws.Range("E3").Value = "2.3 km from center"
$ws.Rows.Item(3).AutoFit()

$ws.Range("A4").Value = "APA Hotel Higashi Shinjuku Kabukicho Tower"
$ws.Range("C4").Value = "Scored 8.1"
$ws.Range("D4").Value = "Very Good`n19,368 reviews"
$ws.Range("E4").Value = "center"
$ws.Rows.Item(4).AutoFit()

$ws.Range("A5").Value = "THE TOURIST HOTEL & Cafe AKIHABARA"
$ws.Range("C5").Value = "Scored 8.7"
$ws.Range("D5").Value = "Excellent`n4,320 reviews"
$ws.Range("E5").Value = "center"
$ws.Rows.Item(5).AutoFit()

$ws.Range("A6").Value = "APA Hotel PRIDE Akasaka Kokkaigijidomae"
$ws.Range("C6").Value = "Scored 8.4"
$ws.Range("D6").Value = "Very Good`n8,157 reviews"
$ws.Range("E6").Value = "center"
$ws.Rows.Item(6).AutoFit()

$ws.Range("A7").Value = "APA Hotel & Resort Ryogoku Eki Tower"
$ws.Range("C7").Value = "Scored 8.2"
$ws.Range("D7").Value = "Very Good`n15,060 reviews"
$ws.Range("E7").Value = "6 km from center"
$ws.Rows.Item(7).AutoFit()

$ws.Range("A8").Value = "Best Western Hotel Fino Tokyo Akihabara"
$ws.Range("C8").Value = "Scored 8.2"
$ws.Range("D8").Value = "Very Good`n3,281 reviews"
$ws.Range("E8").Value = "center"
$ws.Rows.Item(8).AutoFit()

$ws.Range("A9").Value = "Asakusa Ryokan Toukaisou"
$ws.Range("C9").Value = "Scored 8.9"
$ws.Range("D9").Value = "Excellent`n1,038 reviews"
$ws.Range("E9").Value = "center"
$ws.Rows.Item(9).AutoFit()

$ws.Range("A10").Value = "THE KNOT TOKYO Shinjuku"
$ws.Range("C10").Value = "Scored 8.3"
$ws.Range("D10").Value = "Very Good`n10,518 reviews"
$ws.Range("E10").Value = "center"
$ws.Rows.Item(10).AutoFit()

$ws.Range("A11").Value = "Hotel Nihonbashi Saibo"
$ws.Range("C11").Value = "Scored 8.5"
$ws.Range("D11").Value = "Very Good`n2,934 reviews"
$ws.Range("E11").Value = "250 m from center"
$ws.Rows.Item(11).AutoFit()

$ws.Range("A12").Value = "Tabist Urban Stays Asakusa"
$ws.Range("C12").Value = "Scored 8.0"
$ws.Range("D12").Value = "Very Good`n1,642 reviews"
$ws.Range("E12").Value = "350 m from center"
$ws.Rows.Item(12).AutoFit()

$ws.Range("A13").Value = "LANDABOUT TOKYO"
$ws.Range("C13").Value = "Scored 8.7"
$ws.Range("D13").Value = "Excellent`n3,627 reviews"
$ws.Range("E13").Value = "6.1 km from center"
$ws.Rows.Item(13).AutoFit()

$ws.Range("A14").Value = "Restay Frontier (Adult Only)"
$ws.Range("C14").Value = "Scored 8.7"
$ws.Range("D14").Value = "Excellent`n621 reviews"
$ws.Range("E14").Value = "3.3 km from center"
$ws.Rows.Item(14).AutoFit()

$ws.Range("A15").Value = "Akabane Holic Hotel"
$ws.Range("C15").Value = "Scored 8.5"
$ws.Range("D15").Value = "Very Good`n3,967 reviews"
$ws.Range("E15").Value = "8 km from center"
$ws.Rows.Item(15).AutoFit()

$ws.Range("A16").Value = "Toshi Center Hotel"
$ws.Range("C16").Value = "Scored 8.5"
$ws.Range("D16").Value = "Very Good`n3,826 reviews"
$ws.Range("E16").Value = "center"
$ws.Rows.Item(16).AutoFit()

$ws.Range("A17").Value = "HOTEL MYSTAYS Asakusabashi"
$ws.Range("C17").Value = "Scored 8.2"
$ws.Range("D17").Value = "Very Good`n2,403 reviews"
$ws.Range("E17").Value = "300 m from center"
$ws.Rows.Item(17).AutoFit()

$ws.Range("A18").Value = "Forest Hongo by unito"
$ws.Range("C18").Value = "Scored 8.6"
$ws.Range("D18").Value = "Excellent`n1,018 reviews"
$ws.Range("E18").Value = "0.5 km from center"
$ws.Rows.Item(18).AutoFit()

$ws.Range("A19").Value = "SF Heigths"
$ws.Range("C19").Value = "Scored 7.3"
$ws.Range("D19").Value = "Good`n15 reviews"
$ws.Range("E19").Value = "center"
$ws.Rows.Item(19).AutoFit()

$ws.Range("A20").Value = "APA Hotel Nihombashi Bakuroyokoyama Ekimae"
$ws.Range("C20").Value = "Scored 8.4"
$ws.Range("D20").Value = "Very Good`n4,241 reviews"
$ws.Range("E20").Value = "0.6 km from center"
$ws.Rows.Item(20).AutoFit()

$ws.Range("A21").Value = "APA Hotel Asakusa Kuramae Kita"
$ws.Range("C21").Value = "Scored 8.5"
$ws.Range("D21").Value = "Very Good`n5,064 reviews"
$ws.Range("E21").Value = "350 m from center"
$ws.Rows.Item(21).AutoFit()

$ws.Range("A22").Value = "Hotel Asakusa KANNONURA"
$ws.Range("C22").Value = "Scored 8.2"
$ws.Range("D22").Value = "Very Good`n853 reviews"
$ws.Range("E22").Value = "center"
$ws.Rows.Item(22).AutoFit()

$ws.Range("A23").Value = "APA Hotel TKP Nippori Ekimae"
$ws.Range("C23").Value = "Scored 8.1"
$ws.Range("D23").Value = "Very Good`n2,097 reviews"
$ws.Range("E23").Value = "0.9 km from center"
$ws.Rows.Item(23).AutoFit()

$ws.Range("A24").Value = "NEW OPEN!! RUTiLE IKEBUKURO TOKYO"
$ws.Range("C24").Value = "Scored 8.6"
$ws.Range("D24").Value = "Excellent`n20 reviews"
$ws.Range("E24").Value = "4.1 km from center"
$ws.Rows.Item(24).AutoFit()

$ws.Range("A25").Value = "ART HOTELS SHIBUYA"
$ws.Range("C25").Value = "Scored 8.6"
$ws.Range("D25").Value = "Excellent`n1,154 reviews"
$ws.Range("E25").Value = "250 m from center"
$ws.Rows.Item(25).AutoFit()

$ws.Range("A26").Value = "HOTEL TAVINOS Hamamatsucho"
$ws.Range("C26").Value = "Scored 8.2"
$ws.Range("D26").Value = "Very Good`n1,191 reviews"
$ws.Range("E26").Value = "200 m from center"
$ws.Rows.Item(26).AutoFit()
